# Add two new columns (I: "I0", J: "IF") to the worksheet, mirroring the
# existing header style used by columns B-H and filling in the per-row
# numeric data for rows 2-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style/formatting from H1 onto the new header cells I1:J1,
# then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the numeric data for the new columns I and J, rows 2-15.
$data = @{
    2  = @(8, 8)
    3  = @(8, 8)
    4  = @(9, 9)
    5  = @(4, 5)
    6  = @(8, 8)
    7  = @(8, 8)
    8  = @(8, 8)
    9  = @(8, 8)
    10 = @(7, 7)
    11 = @(8, 8)
    12 = @(5, 5)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
